$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "You have received additional news that the disease has sweeped through all neighbouring kingdoms and the situation has gotten severe."
$ws.Range("B16").Value = "Thankfully, the strict border restrictions helped to manage the people entering your kingdom, minimising any risk of potential foreign infection. The additional medicine stock also meant that your kingdom is able to cope should an outbreak occur."
$ws.Range("B18").Value = "Cutting off interactions with neighbouring kingdoms aided in reducing the risk of any infected people from the neighbouring kingdoms entering your own. However, this came at a cost to your kingdom's trade and relationship ties with them."
